$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.994.86'
$ws.Range("E2").Value = '  -1.83%  '

$ws.Range("D3").Value = '1.908.10'
$ws.Range("E3").Value = '  -4.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.65'
$ws.Range("E5").Value = '  -1.30%  '

$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4602'
$ws.Range("E7").Value = '  -1.84%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3830'
$ws.Range("E8").Value = '  -2.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07751'
$ws.Range("E9").Value = '  -2.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9845'
$ws.Range("E10").Value = '  -1.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.08'
$ws.Range("E11").Value = '  -3.18%  '

$ws.Range("D12").Value = '1.888.32'
$ws.Range("E12").Value = '  -7.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.997'
$ws.Range("E13").Value = '  -3.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.688'
$ws.Range("E14").Value = '  -3.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07051'
$ws.Range("E15").Value = '  -1.43%  '

$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.06'
$ws.Range("E17").Value = '  -5.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009551'
$ws.Range("E18").Value = '  -4.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.77'
$ws.Range("E19").Value = '  -3.37%  '

$ws.Range("E20").Value = '  -0.30%  '

$ws.Range("D21").Value = '28.990.75'
$ws.Range("E21").Value = '  -2.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.330'
$ws.Range("E22").Value = '  -4.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.97'
$ws.Range("E23").Value = '  -2.82%  '

$ws.Range("D24").Value = '2.159.41'
$ws.Range("E24").Value = '  -4.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.077'
$ws.Range("E25").Value = '  -2.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.92'
$ws.Range("E26").Value = '  -0.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.23'
$ws.Range("E27").Value = '  -2.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.607'
$ws.Range("E28").Value = '  -6.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.92'
$ws.Range("E29").Value = '  -2.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.830'
$ws.Range("E30").Value = '  -6.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09276'
$ws.Range("E31").Value = '  -1.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8614'
$ws.Range("E32").Value = '  -3.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.116'
$ws.Range("E33").Value = '  -3.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.254'
$ws.Range("E34").Value = '  -7.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.021'
$ws.Range("E35").Value = '  -5.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05738'
$ws.Range("E36").Value = '  -1.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.146'
$ws.Range("E37").Value = '  -2.55%  '

$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02047'
$ws.Range("E39").Value = '  -3.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.477'
$ws.Range("E40").Value = '  -5.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5529'
$ws.Range("E41").Value = '  -4.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1758'
$ws.Range("E42").Value = '  -3.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.286'
$ws.Range("E43").Value = '  -5.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.742'
$ws.Range("E44").Value = '  +3.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5219'
$ws.Range("E45").Value = '  -3.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.26'
$ws.Range("E46").Value = '  -6.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.096'
$ws.Range("E47").Value = '  -2.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06829'
$ws.Range("E48").Value = '  -1.96%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.05'
$ws.Range("E49").Value = '  -2.29%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.782'
$ws.Range("E50").Value = '  -4.76%  '

$ws.Range("B51").Value = 'PEPE'
$ws.Range("C51").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000002570'
$ws.Range("E51").Value = '  -17.42%  '
